$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cell with the default (unstyled) cell format, used to restore
# styling on cells that briefly need a Text number format so that purely
# numeric-looking strings (e.g. "511.50") are not auto-coerced into numbers.
$cleanStyle = $ws.Range("B2").Style

$ws.Range("D2").Value = '56.740.52'
$ws.Range("E2").Value = '  +0.45%  '
$ws.Range("D3").Value = '3.028.61'
$ws.Range("E3").Value = '  +2.56%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '511.50'
$ws.Range("D5").Style = $cleanStyle
$ws.Range("E5").Value = '  +3.11%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '140.45'
$ws.Range("D6").Style = $cleanStyle
$ws.Range("E6").Value = '  +4.90%  '
$ws.Range("E8").Value = '  +2.04%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.13'
$ws.Range("D9").Style = $cleanStyle
$ws.Range("E9").Value = '  +0.06%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.369'
$ws.Range("D11").Style = $cleanStyle
$ws.Range("E11").Value = '  +5.33%  '
$ws.Range("D12").Value = '3.548.19'
$ws.Range("E12").Value = '  +2.50%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.34'
$ws.Range("D14").Style = $cleanStyle
$ws.Range("E14").Value = '  -1.71%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000162'
$ws.Range("D15").Style = $cleanStyle
$ws.Range("E15").Value = '  +4.08%  '
$ws.Range("D16").Value = '56.723.01'
$ws.Range("E16").Value = '  +0.33%  '
$ws.Range("D17").Value = '3.029.50'
$ws.Range("E17").Value = '  +2.18%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '5.93'
$ws.Range("D18").Style = $cleanStyle
$ws.Range("E18").Value = '  -0.40%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.15'
$ws.Range("D19").Style = $cleanStyle
$ws.Range("E19").Value = '  +5.99%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.05'
$ws.Range("D20").Style = $cleanStyle
$ws.Range("E20").Value = '  +4.08%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '333.32'
$ws.Range("D21").Style = $cleanStyle
$ws.Range("E21").Value = '  +5.66%  '
$ws.Range("E22").Value = '  -0.06%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.501'
$ws.Range("D23").Style = $cleanStyle
$ws.Range("E23").Value = '  +3.69%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '64.64'
$ws.Range("D24").Style = $cleanStyle
$ws.Range("E24").Value = '  +3.46%  '
$ws.Range("D25").Value = '3.159.64'
$ws.Range("E25").Value = '  +2.64%  '
$ws.Range("E26").Value = '  +3.05%  '
$ws.Range("E27").Value = '  -0.18%  '
$ws.Range("D28").Value = '0.0₃0929'
$ws.Range("E28").Value = '  +8.49%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.42'
$ws.Range("D29").Style = $cleanStyle
$ws.Range("E29").Value = '  -0.27%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.81'
$ws.Range("D30").Style = $cleanStyle
$ws.Range("E30").Value = '  -2.52%  '
$ws.Range("E31").Value = '  +3.00%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '20.42'
$ws.Range("D32").Style = $cleanStyle
$ws.Range("E32").Value = '  +2.85%  '
$ws.Range("E33").Value = '  +3.18%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '152.68'
$ws.Range("D34").Style = $cleanStyle
$ws.Range("E34").Value = '  -0.01%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.48'
$ws.Range("D35").Style = $cleanStyle
$ws.Range("E35").Value = '  +0.84%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '27.06'
$ws.Range("D36").Style = $cleanStyle
$ws.Range("E36").Value = '  +14.15%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.81'
$ws.Range("D37").Style = $cleanStyle
$ws.Range("E37").Value = '  +2.66%  '
$ws.Range("E38").Value = '  +2.33%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0663'
$ws.Range("D39").Style = $cleanStyle
$ws.Range("E39").Value = '  +2.08%  '
$ws.Range("D40").Value = '3.067.43'
$ws.Range("E40").Value = '  +2.80%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '36.56'
$ws.Range("D41").Style = $cleanStyle
$ws.Range("E41").Value = '  -1.79%  '
$ws.Range("E42").Value = '  -0.04%  '
$ws.Range("E43").Value = '  +3.94%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.661'
$ws.Range("D44").Style = $cleanStyle
$ws.Range("E44").Value = '  +3.69%  '
$ws.Range("D45").Value = '2.213.25'
$ws.Range("E45").Value = '  +3.48%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.35'
$ws.Range("D46").Style = $cleanStyle
$ws.Range("E46").Value = '  +0.86%  '
$ws.Range("E47").Value = '  +6.38%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.933'
$ws.Range("D48").Style = $cleanStyle
$ws.Range("E48").Value = '  +1.99%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '19.81'
$ws.Range("D49").Style = $cleanStyle
$ws.Range("E49").Value = '  +5.00%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '5.84'
$ws.Range("D50").Style = $cleanStyle
$ws.Range("E50").Value = '  +0.08%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0855'
$ws.Range("D51").Style = $cleanStyle
$ws.Range("E51").Value = '  +1.28%  '
